# Add 2022-Q3 data: insert a new quarter sheet after "总计", shifting every
# existing quarter sheet's data down by one position, and append a fresh
# copy of the former last sheet ("2020-Q4") at the end so its unique
# formatting is preserved in the new final slot.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Duplicate the current LAST sheet ("2020-Q4", which carries the
#    workbook's unique header/margin styling) and place the duplicate
#    right after it. That duplicate will become the new, final
#    "2020-Q4" sheet and keeps its data/style untouched.
# ---------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$oldLastSheet = $wb.Worksheets.Item($sheetCount)
$oldLastSheet.Copy($null, $oldLastSheet)

# ---------------------------------------------------------------------
# 2) Duplicate sheet #7 ("2021-Q1", common styling shared by all the
#    other quarter sheets) and place the duplicate right after it. This
#    gives us a ninth, common-styled sheet to become the new "2021-Q1"
#    slot, while the original last sheet (still sitting between them)
#    gets removed since its data lives on in the sheet created in step 1.
# ---------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item(7)
$templateSheet.Copy($null, $templateSheet)

# The original "2020-Q4" sheet (pre-duplication) is now sandwiched at
# position 9 (between the new common-style duplicate and the final
# "2020-Q4" duplicate) - remove it, it is no longer needed.
$spareSheet = $wb.Worksheets.Item(9)
$spareSheet.Delete()

# ---------------------------------------------------------------------
# 3) Write the shifted quarterly data into each slot and rename the
#    tabs into their new positions.
# ---------------------------------------------------------------------

# Slot 2: was "2022-Q2" -> becomes "2022-Q3" (brand-new data)
$s2 = $wb.Worksheets.Item(2)
$s2.Range("D1").Value = "基金规模"
$s2.Range("D2").Value = "0.58"
$s2.Range("E2").Value = "95.06"
$s2.Range("F2").Value = "5.71"
$s2.Range("G2").Value = "0.0331"
$s2.Range("H2").Value = 4
$s2.Name = "2022-Q3"

# Slot 3: was "2022-Q1" -> becomes "2022-Q2" (old 2022-Q2 values)
$s3 = $wb.Worksheets.Item(3)
$s3.Range("D1").Value = "基金规模"
$s3.Range("D2").Value = "0.57"
$s3.Range("E2").Value = "97.64"
$s3.Range("F2").Value = "5.83"
$s3.Range("G2").Value = "0.0332"
$s3.Range("H2").Value = 4
$s3.Name = "2022-Q2"

# Slot 4: was "2021-Q4" -> becomes "2022-Q1" (old 2022-Q1 values)
$s4 = $wb.Worksheets.Item(4)
$s4.Range("D1").Value = "基金规模"
$s4.Range("D2").Value = "0.60"
$s4.Range("E2").Value = "96.69"
$s4.Range("F2").Value = "5.64"
$s4.Range("G2").Value = "0.0338"
$s4.Range("H2").Value = 4
$s4.Name = "2022-Q1"

# Slot 5: was "2021-Q3" -> becomes "2021-Q4" (old 2021-Q4 values;
# label switches from "基金金额" to "基金规模")
$s5 = $wb.Worksheets.Item(5)
$s5.Range("D1").Value = "基金规模"
$s5.Range("D2").Value = "0.60"
$s5.Range("E2").Value = "98.37"
$s5.Range("F2").Value = "6.12"
$s5.Range("G2").Value = "0.0367"
$s5.Range("H2").Value = 3
$s5.Name = "2021-Q4"

# Slot 6: was "2021-Q2" -> becomes "2021-Q3" (old 2021-Q3 values)
$s6 = $wb.Worksheets.Item(6)
$s6.Range("D1").Value = "基金金额"
$s6.Range("D2").Value = "0.54"
$s6.Range("E2").Value = "98.50"
$s6.Range("F2").Value = "5.85"
$s6.Range("G2").Value = "0.0316"
$s6.Range("H2").Value = 4
$s6.Name = "2021-Q3"

# Slot 7: was "2021-Q1" -> becomes "2021-Q2" (old 2021-Q2 values)
$s7 = $wb.Worksheets.Item(7)
$s7.Range("D1").Value = "基金金额"
$s7.Range("D2").Value = "0.56"
$s7.Range("E2").Value = "98.23"
$s7.Range("F2").Value = "6.14"
$s7.Range("G2").Value = "0.0344"
$s7.Range("H2").Value = 4
$s7.Name = "2021-Q2"

# Slot 8: duplicate created in step 2, already holds the original
# "2021-Q1" data/style untouched - just rename the tab.
$s8 = $wb.Worksheets.Item(8)
$s8.Name = "2021-Q1"

# Slot 9: duplicate created in step 1, already holds the original
# "2020-Q4" data/style untouched - just rename the tab.
$s9 = $wb.Worksheets.Item(9)
$s9.Name = "2020-Q4"

# ---------------------------------------------------------------------
# 4) Update the "总计" summary sheet: shift every row down one quarter
#    and append a new row for "2020-Q4".
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.03

$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.03

$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 1
$summary.Range("D4").Value = 0.03

$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 1
$summary.Range("D5").Value = 0.04

$summary.Range("B6").Value = "2021-Q3"
$summary.Range("C6").Value = 1
$summary.Range("D6").Value = 0.03

$summary.Range("B7").Value = "2021-Q2"
$summary.Range("C7").Value = 1
$summary.Range("D7").Value = 0.03

$summary.Range("B8").Value = "2021-Q1"
$summary.Range("C8").Value = 1
$summary.Range("D8").Value = 0.03

$summary.Range("A9").Value = 7
$summary.Range("B9").Value = "2020-Q4"
$summary.Range("C9").Value = 1
$summary.Range("D9").Value = 0.03
